# Repull data, push all data, mean calculation
# Update the dSF (column F) values for several rows to reflect the
# repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -12
$ws.Range("F7").Value = -8
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = 4
$ws.Range("F14").Value = -10
$ws.Range("F19").Value = -1
